# Update the Diagrams deck:
#  1. Fixed "Date" placeholder text 11/6/16 -> 11/7/16 across every slide
#     layout, the slide master, and the notes master.
#  2. Sequence-diagram labels on slide 3: "delete 1" -> "del 1" (both the
#     plain label textbox and the execute("...") textbox), preserving the
#     existing run split caused by the in-place edit.
#  3. Trivial empty <p:timing> scaffold added on slides 3 and 4 (mirrors
#     the one already present on slides 6/7) via an Animation Pane touch.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# 1a. Slide layouts (11 of them)
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $lyt = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $lyt.Shapes "11/7/16"
}

# 1b. Slide master
Set-DatePlaceholderText $master.Shapes "11/7/16"

# 1c. Notes master
Set-DatePlaceholderText $p.NotesMaster.Shapes "11/7/16"

# 2. Slide 3 sequence-diagram edits
$slide3 = $p.Slides.Item(3)

# "delete 1" -> "del 1" (plain label textbox, shape 12 / id 24)
$lbl = $slide3.Shapes.Item(12)
$lbl.TextFrame.TextRange.Characters(1, 7).Text = "del "

# execute("delete 1") -> execute("del 1") (shape 14 / id 26)
$exec = $slide3.Shapes.Item(14)
$exec.TextFrame.TextRange.Characters(10, 7).Text = "del "

# 3. Empty timing scaffolds on slides 3 and 4
foreach ($idx in 3, 4) {
    $slide = $p.Slides.Item($idx)
    $seq = $slide.TimeLine.MainSequence
    $eff = $seq.AddEffect($slide.Shapes.Item(1))
    $eff.Delete()
}
